$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5997
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("M122").Value = -3547
